$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting old rows 13-21 down to 14-22.
$ws.Rows.Item(13).Insert()

# Row 10: Objetivos -> replace stray "Teresa Cristina" value with the real
# (Portuguese) objectives text.
$ws.Range("B10").Value = "Introduzir conceitos teóricos e práticos de Ecotoxicologia Aquática para estudantes de Engenharia Ambiental."
$ws.Range("C10").Value = "Introduzir conceitos teóricos e práticos de Ecotoxicologia Aquática para estudantes de Engenharia Ambiental."

# Row 13 (new row): holds "Docentes responsáveis" value only (no label in A13).
$ws.Range("B13").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C13").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Rows.Item(13).AutoFit()

# Row 14: Programa resumido -> replace stray "Semestral" with real short syllabus (PT).
$ws.Range("B14").Value = "Dar conhecimentos aos alunos de noções básicas teóricas e práticas sobre ecotoxicologia aquática e das técnicas usadas em laboratório para os cultivos e os ensaios com os organismos-teste padronizados."
$ws.Range("C14").Value = "Dar conhecimentos aos alunos de noções básicas teóricas e práticas sobre ecotoxicologia aquática e das técnicas usadas em laboratório para os cultivos e os ensaios com os organismos-teste padronizados."
$ws.Rows.Item(14).RowHeight = 60

# Row 15: Short syllabus label moves here (content unchanged); adjust height to 60.
$ws.Range("A15").Value = "Short syllabus:"
$ws.Rows.Item(15).RowHeight = 60

# Row 16: Programa -> replace stray "01/01/2020" with real full syllabus (PT).
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "A. Teórico: 1. Ecotoxicologia: Introdução, histórico, conceitos; 2. Introdução de agentes químicos no ambiente aquático: biodisponibilidade de contaminantes, efeitos sinérgicos e antagônicos, impactos sobre os sistemas aquáticos; 3.Métodos de ensaios de toxicidade com organismos aquáticos: uso de bioindicadores; B. Prática: 4.Seleção, manutenção e cultivo de organismos aquáticos: boas práticas; 5. 5. Testes de toxicidade com organismos aquáticos."
$ws.Range("C16").Value = "A. Teórico: 1. Ecotoxicologia: Introdução, histórico, conceitos; 2. Introdução de agentes químicos no ambiente aquático: biodisponibilidade de contaminantes, efeitos sinérgicos e antagônicos, impactos sobre os sistemas aquáticos; 3.Métodos de ensaios de toxicidade com organismos aquáticos: uso de bioindicadores; B. Prática: 4.Seleção, manutenção e cultivo de organismos aquáticos: boas práticas; 5. 5. Testes de toxicidade com organismos aquáticos."
$ws.Rows.Item(16).RowHeight = 120

# Row 17: Syllabus label/content moves here (content unchanged); height 120.
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# Row 19: Método label now carries the "Aulas teóricas..." text (unchanged value).
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."
$ws.Range("C19").Value = "Aulas teóricas e práticas. Avaliação baseada em prova, exercício e relatório."

# Row 20: Critério label now carries the "Média ponderada..." text (unchanged value).
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas à prova, exercício e relatório."

# Row 21: Norma de recuperação label now carries the "Nota final..." text; height 60.
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Nota final: NF ≥ 5,0"
$ws.Range("C21").Value = "Nota final: NF ≥ 5,0"
$ws.Rows.Item(21).RowHeight = 60

# Row 22 (new row): Bibliografia label with the new bibliography text; height 120.
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Zagatto, P. A.; Bertoletti, E. Ecotoxicologia aquática– princípios e aplicações. RiMa, 2008.Azevedo, F.A.; Chasin, A.M. As bases toxicológicas da ecotoxicologia. RiMa, 2003.MOZETO, A. A.; UMBUZEIRO, G. A.; JARDIM, W. F. Métodos de coleta, análises físico-químicas e ensaios biológicos e ecotoxicológicos de sedimentos de água doce. São Carlos – SP. Cubo Multimídia & Propaganda, 2006."
$ws.Range("C22").Value = "Zagatto, P. A.; Bertoletti, E. Ecotoxicologia aquática– princípios e aplicações. RiMa, 2008.Azevedo, F.A.; Chasin, A.M. As bases toxicológicas da ecotoxicologia. RiMa, 2003.MOZETO, A. A.; UMBUZEIRO, G. A.; JARDIM, W. F. Métodos de coleta, análises físico-químicas e ensaios biológicos e ecotoxicológicos de sedimentos de água doce. São Carlos – SP. Cubo Multimídia & Propaganda, 2006."
$ws.Rows.Item(22).RowHeight = 120

# Column widths: the "before" workbook had columns 1-2 sharing one <col> span
# (min=1 max=2); the target splits that into separate entries for column 1.
$ws.Columns.Item(1).ColumnWidth = 30.7109375
$ws.Columns.Item(2).ColumnWidth = 60.7109375

Write-Output "done"
